# "Generate Report for Handback" — refresh the handback timestamps for the
# 530290e6-5491-4005-a69b-99d51fd2293c record (row 3 of each sheet) with the
# results of a later handback run. The 5183f9b5... record (row 2) already
# reflects its own handback run and is left untouched.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
# "Latest HO Xliff Generate Date" for the 530290e6 file now reflects the
# de-de handoff timestamp produced by this run.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-06 06:56:44"

# --- zh-cn sheet --------------------------------------------------------
# "Correspond Handoff Datetime" (H) and "Correspond Handback DateTime" (K)
# for the 530290e6 row get the new run's timestamps.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-06 06:56:39"
$wsZhCn.Range("K3").Value = "2016-09-06 06:56:57"

# --- de-de sheet ----------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-09-06 06:56:44"
$wsDeDe.Range("K3").Value = "2016-09-06 06:57:13"
